# regen sval data to filter save games
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @{
    2  = @{ B = 3.182878228561681;  C = 1.65323645889881;  D = 0.7127328510149897; E = 0.4998867070740569; G = 6.048734245549538 }
    3  = @{ B = 0.1554434735375247; C = 0.3375848360084654; D = 3.082599426703578;  E = 0.4998867070740569; G = 4.075514443323626 }
    4  = @{ B = 3.182878228561681;  C = 1.65323645889881;  D = 0.1529057820181812; E = 0.4998867070740569; G = 5.488907176552729 }
    5  = @{ B = 3.182878228561681;  C = 1.65323645889881;  D = 3.082599426703578;  E = 6.48142807727062;   G = 14.40014219143469 }
    6  = @{ B = 3.182878228561681;  C = 1.65323645889881;  D = 0.7127328510149897; E = 0.4998867070740569; G = 6.048734245549538 }
    7  = @{ B = 3.182878228561681;  C = 1.65323645889881;  D = 0.7127328510149897; E = 0.4998867070740569; G = 6.048734245549538 }
    8  = @{ B = 0.7287194209349384; C = 1.65323645889881;  D = 0.1529057820181812; E = 0.4998867070740569; G = 3.034748368925986 }
    9  = @{ B = 1.505614041169197;  C = 1.65323645889881;  D = 0.7127328510149897; E = 0.4998867070740569; G = 4.371470058157054 }
    10 = @{ B = 3.182878228561681;  C = 1.65323645889881;  D = 0.7127328510149897; E = 0.4998867070740569; G = 6.048734245549538 }
    11 = @{ B = 3.182878228561681;  C = 1.65323645889881;  D = 0.7127328510149897; E = 0.4998867070740569; G = 6.048734245549538 }
    12 = @{ B = 0.3464964993005633; C = 0.3375848360084654; D = 3.082599426703578;  E = 0.4998867070740569; G = 4.266567469086664 }
}

foreach ($row in $data.Keys) {
    $vals = $data[$row]
    foreach ($col in $vals.Keys) {
        $ws.Range("$col$row").Value = $vals[$col]
    }
}
